$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
    $ws.Range($cellRef).Style = "Normal"
}

Set-CellText 'D2' '29.080.20'
Set-CellText 'E2' '  -0.02%  '
Set-CellText 'D3' '1.832.95'
Set-CellText 'E3' '  -0.20%  '
Set-CellText 'D4' '1.003'
Set-CellText 'E4' '  +0.15%  '
Set-CellText 'D5' '243.71'
Set-CellText 'E5' '  +0.44%  '
Set-CellText 'D6' '0.6281'
Set-CellText 'E6' '  +0.33%  '
Set-CellText 'E7' '  +0.04%  '
Set-CellText 'D8' '0.07464'
Set-CellText 'E8' '  -1.20%  '
Set-CellText 'D9' '0.2920'
Set-CellText 'E9' '  +0.01%  '
Set-CellText 'E10' '  +1.90%  '
Set-CellText 'E11' '  -0.39%  '
Set-CellText 'D12' '1.835.20'
Set-CellText 'E12' '  -0.15%  '
Set-CellText 'D13' '4.979'
Set-CellText 'E13' '  +0.51%  '
Set-CellText 'D14' '0.6682'
Set-CellText 'E14' '  +0.54%  '
Set-CellText 'D15' '82.65'
Set-CellText 'E15' '  -0.15%  '
Set-CellText 'D16' '0.000009318'
Set-CellText 'E16' '  -6.14%  '
Set-CellText 'E17' '  -0.07%  '
Set-CellText 'D18' '29.085.36'
Set-CellText 'E18' '  -0.12%  '
Set-CellText 'D19' '2.074.81'
Set-CellText 'E19' '  -0.57%  '
Set-CellText 'D20' '12.58'
Set-CellText 'E20' '  +1.86%  '
Set-CellText 'D21' '222.85'
Set-CellText 'E21' '  -1.28%  '
Set-CellText 'E22' '  +0.17%  '
Set-CellText 'D23' '7.121'
Set-CellText 'E23' '  -0.85%  '
Set-CellText 'D24' '1.004'
Set-CellText 'E24' '  +0.10%  '
Set-CellText 'D25' '160.22'
Set-CellText 'E25' '  +0.91%  '
Set-CellText 'D26' '0.1396'
Set-CellText 'E26' '  +2.15%  '
Set-CellText 'E27' '  +0.34%  '
Set-CellText 'D28' '17.90'
Set-CellText 'E28' '  -0.02%  '
Set-CellText 'D29' '1.501'
Set-CellText 'E29' '  +0.73%  '
Set-CellText 'D30' '0.05771'
Set-CellText 'E30' '  +10.89%  '
Set-CellText 'D31' '4.149'
Set-CellText 'E31' '  +1.73%  '
Set-CellText 'D32' '4.066'
Set-CellText 'E32' '  +0.79%  '
Set-CellText 'D33' '1.205'
Set-CellText 'E33' '  +0.36%  '
Set-CellText 'D34' '0.7499'
Set-CellText 'E34' '  +1.52%  '
Set-CellText 'D35' '1.849'
Set-CellText 'E35' '  -0.43%  '
Set-CellText 'E36' '  -0.48%  '
Set-CellText 'D37' '2.622'
Set-CellText 'E37' '  -3.09%  '
Set-CellText 'D38' '1.227.81'
Set-CellText 'E38' '  -1.84%  '
Set-CellText 'D39' '2.754'
Set-CellText 'E39' '  -0.42%  '
Set-CellText 'D40' '0.01784'
Set-CellText 'E40' '  +0.04%  '
Set-CellText 'D41' '6.552'
Set-CellText 'E41' '  +2.86%  '
Set-CellText 'D42' '0.8921'
Set-CellText 'E42' '  -0.07%  '
Set-CellText 'D44' '101.83'
Set-CellText 'E44' '  +0.17%  '
Set-CellText 'D45' '1.987.52'
Set-CellText 'E45' '  +0.38%  '
Set-CellText 'B46' 'BabyDogeCoin'
Set-CellText 'C46' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-CellText 'D46' '0.00000000125'
Set-CellText 'E46' '  +0.06%  '
Set-CellText 'B47' 'Aave'
Set-CellText 'C47' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText 'D47' '65.68'
Set-CellText 'E47' '  +2.53%  '
Set-CellText 'D48' '0.07781'
Set-CellText 'E48' '  +16.47%  '
Set-CellText 'D49' '0.5101'
Set-CellText 'E49' '  -0.39%  '
Set-CellText 'E50' '  +0.99%  '
Set-CellText 'D51' '9.015'
Set-CellText 'E51' '  +1.93%  '
